$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "year" value per group of 10 rows (2014..2023), each
# currently stored as a text date string ("31/12/2014", ...). Replace them
# with plain numeric years.
$startYear = 2014
for ($block = 0; $block -lt 10; $block++) {
    $year = $startYear + $block
    $rowStart = 2 + ($block * 10)
    $rowEnd = $rowStart + 9
    $ws.Range("A$rowStart`:A$rowEnd").Value = $year
}

# Remove the selection on the sheet view (tidy-up left over from editing),
# restoring it to the default (A1 selected).
$ws.Range("A1").Select()
